# Apply 2024-10-18 data update to violent-crime-full-year workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("K2").Value = 5827
$ws.Range("K3").Value = 5994
$ws.Range("K4").Value = 1248
$ws.Range("K5").Value = 425
$ws.Range("K6").Value = 6601
$ws.Range("K7").Value = 20095

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("K4").Value = 77
$ws.Range("K5").Value = 50
$ws.Range("K8").Value = 1329
$ws.Range("K9").Value = 87
$ws.Range("K10").Value = 113
$ws.Range("K18").Value = 132
$ws.Range("K22").Value = 57
$ws.Range("K24").Value = 59
$ws.Range("K25").Value = 95
$ws.Range("J27").Value = 181
$ws.Range("K27").Value = 187
$ws.Range("K29").Value = 1099
$ws.Range("K31").Value = 222
$ws.Range("K37").Value = 676
$ws.Range("K41").Value = 138
$ws.Range("K42").Value = 744
$ws.Range("K43").Value = 172
$ws.Range("K44").Value = 171
$ws.Range("K47").Value = 138
$ws.Range("K52").Value = 530
$ws.Range("K53").Value = 255
$ws.Range("K54").Value = 390
$ws.Range("K59").Value = 33
$ws.Range("K60").Value = 121
$ws.Range("J63").Value = 114
$ws.Range("K63").Value = 61
$ws.Range("K64").Value = 128
$ws.Range("K65").Value = 466
$ws.Range("K67").Value = 783
$ws.Range("K74").Value = 21
$ws.Range("K76").Value = 272
$ws.Range("K77").Value = 140
$ws.Range("K78").Value = 226
$ws.Range("K83").Value = 448
$ws.Range("K84").Value = 159
$ws.Range("K85").Value = 944
$ws.Range("K87").Value = 34
$ws.Range("K88").Value = 215
$ws.Range("K89").Value = 295
$ws.Range("K90").Value = 186
$ws.Range("K91").Value = 228
$ws.Range("K95").Value = 336
$ws.Range("K96").Value = 212
$ws.Range("K97").Value = 160
$ws.Range("K101").Value = 20095

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("K6").Value = 91
$ws.Range("K7").Value = 212

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("K4").Value = 34
$ws.Range("K6").Value = 89
$ws.Range("K7").Value = 295

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("K2").Value = 308
$ws.Range("K3").Value = 322
$ws.Range("K5").Value = 28
$ws.Range("K7").Value = 944

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("K4").Value = 29
$ws.Range("K6").Value = 187
$ws.Range("K7").Value = 530

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("K3").Value = 67
$ws.Range("K6").Value = 110
$ws.Range("K7").Value = 255

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("K2").Value = 369
$ws.Range("K3").Value = 405
$ws.Range("K6").Value = 444
$ws.Range("K7").Value = 1329

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("K2").Value = 155
$ws.Range("K3").Value = 160
$ws.Range("K7").Value = 448

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("K3").Value = 318
$ws.Range("K6").Value = 257

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("K2").Value = 114
$ws.Range("K7").Value = 336

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("K6").Value = 196
$ws.Range("K7").Value = 676

$ws = $wb.Worksheets.Item('New City')
$ws.Range("K3").Value = 113
$ws.Range("K7").Value = 466

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("K4").Value = 10
$ws.Range("K7").Value = 222

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("K2").Value = 220
$ws.Range("K6").Value = 220
$ws.Range("K7").Value = 783

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("K2").Value = 55
$ws.Range("K7").Value = 159

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("K2").Value = 64
$ws.Range("K6").Value = 209
$ws.Range("K7").Value = 390

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("K3").Value = 393
$ws.Range("K6").Value = 313
$ws.Range("K7").Value = 1099

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("K2").Value = 45
$ws.Range("K7").Value = 171

$ws = $wb.Worksheets.Item('River North')
$ws.Range("K3").Value = 52
$ws.Range("K7").Value = 272

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range("K2").Value = 49
$ws.Range("K7").Value = 138

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("K3").Value = 229
$ws.Range("K6").Value = 278
$ws.Range("K7").Value = 744

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("K6").Value = 52
$ws.Range("K7").Value = 113

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("K6").Value = 80
$ws.Range("K7").Value = 226

$ws = $wb.Worksheets.Item('Dunning')
$ws.Range("K2").Value = 22
$ws.Range("K7").Value = 59

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("K3").Value = 110
$ws.Range("K7").Value = 228

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("K3").Value = 36
$ws.Range("K7").Value = 128

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("K4").Value = 16
$ws.Range("K7").Value = 132

$ws = $wb.Worksheets.Item('East Side')
$ws.Range("K2").Value = 34
$ws.Range("K3").Value = 33
$ws.Range("K7").Value = 95

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("K2").Value = 39
$ws.Range("K7").Value = 138

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range("K3").Value = 30
$ws.Range("K7").Value = 87

$ws = $wb.Worksheets.Item('Montclare')
$ws.Range("K6").Value = 10
$ws.Range("K7").Value = 33

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("K2").Value = 31
$ws.Range("K7").Value = 160

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("K3").Value = 64
$ws.Range("K7").Value = 215

$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range("K3").Value = 14
$ws.Range("K7").Value = 50

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("J4").Value = 23
$ws.Range("K4").Value = 22
$ws.Range("K6").Value = 68
$ws.Range("J7").Value = 181
$ws.Range("K7").Value = 187

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("K2").Value = 68
$ws.Range("K6").Value = 43
$ws.Range("K7").Value = 186

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range("K6").Value = 36
$ws.Range("K7").Value = 121

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("K4").Value = 23
$ws.Range("K7").Value = 172

$ws = $wb.Worksheets.Item('Clearing')
$ws.Range("K2").Value = 27
$ws.Range("K7").Value = 57

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range("K2").Value = 59
$ws.Range("K7").Value = 140

$ws = $wb.Worksheets.Item('Archer Heights')
$ws.Range("K3").Value = 17
$ws.Range("K7").Value = 77

$ws = $wb.Worksheets.Item('Ukrainian Village')
$ws.Range("K6").Value = 18
$ws.Range("K7").Value = 34

$ws = $wb.Worksheets.Item('Printers Row')
$ws.Range("K6").Value = 12
$ws.Range("K7").Value = 21
